# The <id>...</id> tag for this entry was previously split across three
# runs: "<id>" (Courier New / color 7f6000 / sz 18), "p167r_1" (plain
# run), and "</id>" (Courier New / color 7f6000 / sz 18). Collapse them
# into a single run "<id>p167r_1</id>" carrying the Courier New/7f6000/18
# formatting, as newly downloaded tc/tcn/tl content now ships as one run.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p167r_1</id>",  # Find What
    $true,                # MatchCase
    $false,               # MatchWholeWord
    $false,               # MatchWildcards
    $false,               # MatchSoundsLike
    $false,               # MatchAllWordForms
    $true,                # Forward
    1,                    # Wrap (wdFindContinue)
    $false,               # Format
    "<id>p167r_1</id>",  # Replace With
    2                     # Replace (wdReplaceAll)
)
